$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: Breast ---
$ws.Range("B3").Value = 393
$ws.Range("C3").Value = 31
$ws.Range("E3").Value = "noch nicht mit GMD durchgelaufen"

# --- Row 4: Diabetes ---
$ws.Range("B4").Value = 768
$ws.Range("C4").Value = 9

# --- Row 8: Madelon ---
$ws.Range("B8").Value = 330
$ws.Range("C8").Value = 501
$ws.Range("E8").Value = "noch nicht mit GMD durchgelaufen"

# --- Fill column D formula (C-1) for rows 3..10 as one shared formula ---
$ws.Range("D3:D10").Formula = "=C3-1"

# --- column E width (closest achievable to 32.140625 in this runtime) ---
$ws.Columns.Item(5).ColumnWidth = 31.3

# --- selection ---
$ws.Range("D9").Select() | Out-Null
